# Update leve-profit calculation columns (H-N) across the ALC/ARM/BSM/CRP/
# GSM/LTW/WVR sheets to reflect refreshed market-board pricing data.
# Each block re-points the active worksheet and rewrites the affected
# price/profit cells for specific leve rows.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 550
$ws.Range("I12").Value = 100
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 100
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 70
$ws.Range("N12").Value = -1340
$ws.Range("H87").Value = 18786.889
$ws.Range("J87").Value = 18786.889
$ws.Range("L87").Value = 18786.889
$ws.Range("N87").Value = -21282.889
$ws.Range("H90").Value = 18786.889
$ws.Range("J90").Value = 18786.889
$ws.Range("L90").Value = 56360.667
$ws.Range("N90").Value = -68840.667
$ws.Range("H98").Value = 40239.8
$ws.Range("I98").Value = 1450.8823
$ws.Range("J98").Value = 260043.67
$ws.Range("K98").Value = 1450.8823
$ws.Range("L98").Value = 260043.67
$ws.Range("M98").Value = 47.11770000000001
$ws.Range("N98").Value = -263039.67
$ws.Range("H101").Value = 660.2
$ws.Range("I101").Value = 660.2
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 1980.6
$ws.Range("L101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -358.6000000000001
$ws.Range("H122").Value = 40239.8
$ws.Range("I122").Value = 1450.8823
$ws.Range("J122").Value = 260043.67
$ws.Range("K122").Value = 4352.6469
$ws.Range("L122").Value = 780131.01
$ws.Range("M122").Value = -1902.6469
$ws.Range("N122").Value = -785031.01
$ws.Range("H138").Value = 1788.7391
$ws.Range("I138").Value = 1047.898
$ws.Range("J138").Value = 2632.9534
$ws.Range("K138").Value = 3143.694
$ws.Range("L138").Value = 7898.860199999999
$ws.Range("M138").Value = 1996.306
$ws.Range("N138").Value = -18178.8602

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1486.2245
$ws.Range("I74").Value = 1333.186
$ws.Range("K74").Value = 1333.186
$ws.Range("M74").Value = -459.1859999999999
$ws.Range("H77").Value = 1486.2245
$ws.Range("I77").Value = 1333.186
$ws.Range("K77").Value = 6665.929999999999
$ws.Range("M77").Value = -2297.929999999999
$ws.Range("H80").Value = 38196.2
$ws.Range("J80").Value = 38196.2
$ws.Range("L80").Value = 38196.2
$ws.Range("N80").Value = -40192.2
$ws.Range("H83").Value = 38196.2
$ws.Range("J83").Value = 38196.2
$ws.Range("L83").Value = 114588.6
$ws.Range("N83").Value = -124572.6
$ws.Range("H96").Value = 31450
$ws.Range("J96").Value = 31450
$ws.Range("L96").Value = 31450
$ws.Range("N96").Value = -36942
$ws.Range("H105").Value = 49181
$ws.Range("J105").Value = 49181
$ws.Range("L105").Value = 49181
$ws.Range("N105").Value = -56169
$ws.Range("H110").Value = 1603.6364
$ws.Range("I110").Value = 1642.5807
$ws.Range("K110").Value = 1642.5807
$ws.Range("M110").Value = 402.4193

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1936.625
$ws.Range("I86").Value = 1956.1428
$ws.Range("J86").Value = 1800
$ws.Range("K86").Value = 1956.1428
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -833.1428000000001
$ws.Range("N86").Value = -4046
$ws.Range("H89").Value = 1936.625
$ws.Range("I89").Value = 1956.1428
$ws.Range("J89").Value = 1800
$ws.Range("K89").Value = 9780.714
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -4164.714
$ws.Range("N89").Value = -20232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 47826.75
$ws.Range("J81").Value = 47826.75
$ws.Range("L81").Value = 47826.75
$ws.Range("N81").Value = -49822.75
$ws.Range("H84").Value = 47826.75
$ws.Range("J84").Value = 47826.75
$ws.Range("L84").Value = 143480.25
$ws.Range("N84").Value = -153464.25
$ws.Range("H88").Value = 35405.332
$ws.Range("J88").Value = 35405.332
$ws.Range("L88").Value = 35405.332
$ws.Range("N88").Value = -36217.332
$ws.Range("H91").Value = 35405.332
$ws.Range("J91").Value = 35405.332
$ws.Range("L91").Value = 35405.332
$ws.Range("N91").Value = -38213.332
$ws.Range("H122").Value = 71519.53
$ws.Range("I122").Value = 93236
$ws.Range("J122").Value = 941
$ws.Range("K122").Value = 279708
$ws.Range("L122").Value = 2823
$ws.Range("M122").Value = -277258
$ws.Range("N122").Value = -7723
$ws.Range("H124").Value = 33728.6
$ws.Range("J124").Value = 33728.6
$ws.Range("L124").Value = 33728.6
$ws.Range("N124").Value = -38638.6
$ws.Range("H125").Value = 33997.332
$ws.Range("J125").Value = 33997.332
$ws.Range("L125").Value = 33997.332
$ws.Range("N125").Value = -38917.332
$ws.Range("H134").Value = 610220.2
$ws.Range("I134").Value = 1218.2667
$ws.Range("J134").Value = 1752098.8
$ws.Range("K134").Value = 3654.800099999999
$ws.Range("L134").Value = 5256296.4
$ws.Range("M134").Value = -1119.800099999999
$ws.Range("N134").Value = -5261366.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 42131
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 42131
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H86").Value = 30135
$ws.Range("J86").Value = 30135
$ws.Range("L86").Value = 30135
$ws.Range("N86").Value = -32507
$ws.Range("H88").Value = 48575
$ws.Range("J88").Value = 48575
$ws.Range("L88").Value = 48575
$ws.Range("N88").Value = -49477
$ws.Range("H89").Value = 30135
$ws.Range("J89").Value = 30135
$ws.Range("L89").Value = 90405
$ws.Range("N89").Value = -102261
$ws.Range("H91").Value = 48575
$ws.Range("J91").Value = 48575
$ws.Range("L91").Value = 48575
$ws.Range("N91").Value = -51695
$ws.Range("H127").Value = 46654.668
$ws.Range("J127").Value = 46654.668
$ws.Range("L127").Value = 46654.668
$ws.Range("N127").Value = -56574.668
$ws.Range("H131").Value = 29330.666
$ws.Range("J131").Value = 29330.666
$ws.Range("L131").Value = 29330.666
$ws.Range("N131").Value = -39410.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2545.6428
$ws.Range("I7").Value = 1922.5238
$ws.Range("K7").Value = 1922.5238
$ws.Range("M7").Value = -1810.5238
$ws.Range("H55").Value = 478.86667
$ws.Range("I55").Value = 346.75
$ws.Range("J55").Value = 629.8570999999999
$ws.Range("K55").Value = 346.75
$ws.Range("L55").Value = 629.8570999999999
$ws.Range("M55").Value = -173.75
$ws.Range("N55").Value = -975.8570999999999
$ws.Range("H88").Value = 43178.332
$ws.Range("J88").Value = 43178.332
$ws.Range("L88").Value = 43178.332
$ws.Range("N88").Value = -44034.332
$ws.Range("H91").Value = 43178.332
$ws.Range("J91").Value = 43178.332
$ws.Range("L91").Value = 43178.332
$ws.Range("N91").Value = -46142.332
$ws.Range("H126").Value = 2545.6428
$ws.Range("I126").Value = 1922.5238
$ws.Range("K126").Value = 5767.5714
$ws.Range("M126").Value = -3297.5714
$ws.Range("H131").Value = 33254
$ws.Range("J131").Value = 33254
$ws.Range("L131").Value = 33254
$ws.Range("N131").Value = -43334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 961.8125
$ws.Range("I126").Value = 1013.1429
$ws.Range("J126").Value = 602.5
$ws.Range("K126").Value = 3039.4287
$ws.Range("L126").Value = 1807.5
$ws.Range("M126").Value = -569.4287000000004
$ws.Range("N126").Value = -6747.5
